$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New version of main: rolling input data now runs through 2019-08-12 ---
# Column A = date, column B = daily log note (shared string).
$dates = @{}
$notes = @{}

$dates[2] = 43605
$dates[3] = 43606
$dates[4] = 43607
$dates[5] = 43608
$dates[6] = 43609
$dates[7] = 43612
$dates[8] = 43613
$dates[9] = 43614
$dates[10] = 43615
$dates[11] = 43616
$dates[12] = 43619
$dates[13] = 43620
$dates[14] = 43621
$dates[15] = 43622
$dates[16] = 43623
$dates[17] = 43626
$dates[18] = 43627
$dates[19] = 43628
$dates[20] = 43629
$dates[21] = 43630
$dates[22] = 43633
$dates[23] = 43634
$dates[24] = 43635
$dates[25] = 43636
$dates[26] = 43637
$dates[27] = 43640
$dates[28] = 43641
$dates[29] = 43642
$dates[30] = 43643
$dates[31] = 43644
$dates[32] = 43647
$dates[33] = 43648
$dates[34] = 43649
$dates[35] = 43650
$dates[36] = 43651
$dates[37] = 43652
$dates[38] = 43653
$dates[39] = 43654
$dates[40] = 43655
$dates[41] = 43656
$dates[42] = 43657
$dates[43] = 43658
$dates[44] = 43659
$dates[45] = 43660
$dates[46] = 43661
$dates[47] = 43662
$dates[48] = 43663
$dates[49] = 43664
$dates[50] = 43665
$dates[51] = 43666
$dates[52] = 43667
$dates[53] = 43668
$dates[54] = 43669
$dates[55] = 43670
$dates[56] = 43671
$dates[57] = 43672
$dates[58] = 43673
$dates[59] = 43674
$dates[60] = 43675
$dates[61] = 43676
$dates[62] = 43677
$dates[63] = 43678
$dates[64] = 43679
$dates[65] = 43680
$dates[66] = 43681
$dates[67] = 43682
$dates[68] = 43683
$dates[69] = 43684
$dates[70] = 43685
$dates[71] = 43686
$dates[72] = 43687
$dates[73] = 43688
$dates[74] = 43689

$notes[2] = 'First day'
$notes[3] = 'Onboarding, setting up accounts, etc'
$notes[4] = 'Computer locked out, worked from uni on collecting lit'
$notes[5] = 'Onboarding, setting up accounts, etc'
$notes[6] = 'Setting up writing document'
$notes[7] = 'Started working on core functionality, program takes in a list of authors and tweets, uses the author and @s to build and edge list, feeds the edgelist into the networkx package, generates netwrok statsitics for each user in the network and then reports the top 10'
$notes[8] = 'Work on core program'
$notes[9] = 'Work on core program, added basic graphing'
$notes[10] = 'Work on core program, wrapped in twitter bot for demonstration'
$notes[11] = 'Feedback on core program that they are intrested in an input user rather than top 10, and also that they can''t show nonvalidated users in the results, not intrested in the twitter bot functionality'
$notes[12] = 'Implimented hashing function to hide all accounts other than the valid user'
$notes[13] = 'Reworked main program to authenticate from Twitter to get more data'
$notes[14] = 'Investigated methods for identifying valid users'
$notes[15] = 'Implmented basic method for getting valid users from tweet objects, can get senders, but not recivers'
$notes[16] = 'Worked on setting up write up doctument'
$notes[17] = 'Implimented getting recivers valid status by making extra API calls'
$notes[18] = 'Recivers can be identified as valid if the tweet is a retweet - added that functionality and now only have to make 20% of API calls otherwise would'
$notes[19] = 'Adapted hasing so valid users and the given user is shown, built verified user generator which will create a JSON of all verified users so no API calls need to be made by the other scripts'
$notes[20] = 'Pivoted main program so it displayed the given user in a distribution of valid users based on percentile rank within the network metrics'
$notes[21] = 'Work on literature'
$notes[22] = 'Research on better graphing tools'
$notes[23] = 'Added a better graphing tool which is written in R, python saves data and then calls an R script which takes the data and saves an HTML graph, finally python opens the HTML to show the graph'
$notes[24] = 'AL'
$notes[25] = 'Work on verified user generator - will take almost 4 days to run, but should only need run once'
$notes[26] = 'Revewing literature'
$notes[27] = 'Built topic modeller script - this takes a given user, scans their sent tweets and returns their most used hashtags and topcis, based on  topic modelling'
$notes[28] = 'Added sentiment analyis to main - it can now report the sentiment of the whole network as well as tweets about the given user'
$notes[29] = 'Changed mian to use full text by using the extended option in tweepy - but retweet text is in a different place from normal mentions so added if statmenets'
$notes[30] = 'AL'
$notes[31] = 'Revewing literature'
$notes[32] = 'Writing, Looking at previous dissertations and planning the structure of mine, added cheating so tweets send by or about the given user can be added to the network artificially (increases bias)'
$notes[33] = 'Bugfixing, using full retweeted text (which is the only way to avoid truncated text) cuts off the ''@'' which stops retweets being used in the network. Have now manually concatinating the frist mention (in a retweet first mention is always the trageted user) back onto each text object'
$notes[34] = 'Made valid user generator robust to dissconnection and changed it to output a dictonary of ID:screen_name so that in future fast searches can be performed on the ID which is faster to get. Added time feature to topic modeller so you can ask for tweets from X days ago only. Started building sentiment checker.py which will scan for mentions of a user (with the option of adding a topic) and will return the sentiment score of those tweets'
$notes[35] = 'Fixing disconnection bug in Verified_users_generator, building sentiment checker.py'

foreach ($r in $dates.Keys) {
    $ws.Cells.Item($r, 1).Value = $dates[$r]
}

foreach ($r in $notes.Keys) {
    $ws.Cells.Item($r, 2).Value = $notes[$r]
}

# Re-apply the existing built-in date format (numFmtId 14, style index 1)
# to every date cell, including the newly appended rows, by copying the
# format from an already-styled date cell instead of creating a new
# custom number format.
$src = $ws.Range("A2")
[void]$src.Copy()
$dst = $ws.Range("A2:A74")
[void]$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false

[void]$ws.Range("B25").Select()
